$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying material-line data was reordered: row 3's ink/varnish detail
# swapped places with row 4's, and row 6's plate detail swapped places with
# row 7's. Column A (and any other columns that already hold identical text
# on both rows) are unaffected. We read the current values first, then write
# them back swapped.
#
# Columns E (quantity, e.g. "0.08"/"0.14"/"1.00"/"2.00") look like numbers,
# so a plain .Value assignment would be auto-converted to a numeric type.
# To keep them stored as literal text (matching the original workbook),
# those columns are written with a leading apostrophe (quote-prefix), which
# is how Excel itself preserves numeric-looking text.

function Swap-Cell([string]$addrA, [string]$addrB, [bool]$asText) {
    $valA = $ws.Range($addrA).Value2
    $valB = $ws.Range($addrB).Value2
    if ($asText) {
        $ws.Range($addrA).Value = "'" + $valB
        $ws.Range($addrB).Value = "'" + $valA
    } else {
        $ws.Range($addrA).Value = $valB
        $ws.Range($addrB).Value = $valA
    }
}

# Row 3 <-> Row 4
Swap-Cell "D3" "D4" $false
Swap-Cell "E3" "E4" $true
Swap-Cell "H3" "H4" $false

# Row 6 <-> Row 7
Swap-Cell "B6" "B7" $false
Swap-Cell "C6" "C7" $false
Swap-Cell "D6" "D7" $false
Swap-Cell "E6" "E7" $true
Swap-Cell "H6" "H7" $false
